$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new weekly record as row 28, pushing the existing
# rows 28..94 down to 29..95 (dimension grows from R94 to R95).
$ws.Rows.Item(28).Insert()

$ws.Cells.Item(28, 1).Value = 9
$ws.Cells.Item(28, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(28, 3).Value = "Metropolitana"
$ws.Cells.Item(28, 4).Value = 44662
$ws.Cells.Item(28, 5).Value = 13
$ws.Cells.Item(28, 6).Value = 100114007
$ws.Cells.Item(28, 7).Value = "Jengibre"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 610
$ws.Cells.Item(28, 11).Value = 12000
$ws.Cells.Item(28, 12).Value = 13000
$ws.Cells.Item(28, 13).Value = 12500
$ws.Cells.Item(28, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(28, 15).Value = "Perú"
$ws.Cells.Item(28, 16).Value = 962
$ws.Cells.Item(28, 17).Value = 13
$ws.Cells.Item(28, 18).Value = "Hortaliza"
